# Add links to the GitHub Classroom assignments in column E (Assignment)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

$ws.Range("E2").Value = "[Git](https://classroom.github.com/a/MkPXmUFj)"
$ws.Range("E3").Value = "[Machine learning](https://classroom.github.com/a/W8Y8nIDh)"
$ws.Range("E4").Value = "[Iterations](https://classroom.github.com/a/axYOETtq)"
$ws.Range("E5").Value = "[Parallel computing](https://classroom.github.com/a/_w48Vu_p)"
$ws.Range("E6").Value = "[Many Models](https://classroom.github.com/a/oNJOK_ws)/Making maps"

# Update the selected cell to match the saved view state
$ws.Range("F11").Select()
